$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The participant "Anna" had been lost; re-add her. This pushes the existing
# rows 9-15 entries down by one row (both value and formatting move - i.e.
# a genuine row shift), and the newly vacated row 9 becomes Anna's row.
#
# Work bottom-up, cell by cell (not as one big multi-column range) so that
# only the cells that actually hold data are touched/moved - mirrors exactly
# what the source row contained, leaving any other cell untouched.

# row 15 -> row 16
$ws.Range("B15").Cut($ws.Range("B16")) | Out-Null
$ws.Range("C15").Cut($ws.Range("C16")) | Out-Null
$ws.Range("D15").Cut($ws.Range("D16")) | Out-Null
$ws.Range("E15").Cut($ws.Range("E16")) | Out-Null
$ws.Range("F15").Cut($ws.Range("F16")) | Out-Null

# row 14 -> row 15
$ws.Range("B14").Cut($ws.Range("B15")) | Out-Null
$ws.Range("C14").Cut($ws.Range("C15")) | Out-Null
$ws.Range("F14").Cut($ws.Range("F15")) | Out-Null

# row 13 -> row 14
$ws.Range("B13").Cut($ws.Range("B14")) | Out-Null
$ws.Range("C13").Cut($ws.Range("C14")) | Out-Null
$ws.Range("E13").Cut($ws.Range("E14")) | Out-Null
$ws.Range("F13").Cut($ws.Range("F14")) | Out-Null

# row 12 -> row 13
$ws.Range("B12").Cut($ws.Range("B13")) | Out-Null
$ws.Range("C12").Cut($ws.Range("C13")) | Out-Null
$ws.Range("E12").Cut($ws.Range("E13")) | Out-Null
$ws.Range("F12").Cut($ws.Range("F13")) | Out-Null
$ws.Range("M12").Cut($ws.Range("M13")) | Out-Null

# row 11 -> row 12
$ws.Range("B11").Cut($ws.Range("B12")) | Out-Null
$ws.Range("C11").Cut($ws.Range("C12")) | Out-Null
$ws.Range("D11").Cut($ws.Range("D12")) | Out-Null
$ws.Range("E11").Cut($ws.Range("E12")) | Out-Null
$ws.Range("F11").Cut($ws.Range("F12")) | Out-Null

# row 10 -> row 11
$ws.Range("B10").Cut($ws.Range("B11")) | Out-Null
$ws.Range("C10").Cut($ws.Range("C11")) | Out-Null
$ws.Range("D10").Cut($ws.Range("D11")) | Out-Null
$ws.Range("E10").Cut($ws.Range("E11")) | Out-Null
$ws.Range("F10").Cut($ws.Range("F11")) | Out-Null

# row 9 -> row 10
$ws.Range("B9").Cut($ws.Range("B10")) | Out-Null
$ws.Range("C9").Cut($ws.Range("C10")) | Out-Null
$ws.Range("D9").Cut($ws.Range("D10")) | Out-Null
$ws.Range("E9").Cut($ws.Range("E10")) | Out-Null
$ws.Range("F9").Cut($ws.Range("F10")) | Out-Null

# New row 9: Anna
$ws.Range("B9").Value2 = "Anna"
$ws.Range("E9").Value2 = 0.52083333333333337

# A couple of cells have no counterpart in the row feeding them, so after the
# shift they are genuinely empty (not just "cut-vacated") - fully clear them.
$ws.Range("M12").Clear() | Out-Null
$ws.Range("D15").Clear() | Out-Null
$ws.Range("E15").Clear() | Out-Null

Write-Output "done"
